# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-column-name suffixes to the concrete
#   format versions they actually represent (FV2310 / FV2404).
# - Turn the sheet's data range into a proper Excel Table (Table1).
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -----------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newFV2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
$newFV2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt 10; $i++) {
    # Columns A-J (1-10) carry the "_old" -> "_FV2310" headers
    $ws.Cells.Item(1, $i + 1).Value = $newFV2310Headers[$i]
    # Columns L-U (12-21) carry the "_new" -> "_FV2404" headers
    $ws.Cells.Item(1, $i + 12).Value = $newFV2404Headers[$i]
}
# Column K (11) stays "diff" - untouched.

# --- 2. Turn the used range into an Excel Table ---------------------------
$tableRange = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Headers renamed, Table1 created, header row frozen."
